$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Merge the "The " / "software was created..." runs in the
#    "Technical description" paragraph into a single run.
# ------------------------------------------------------------------
$d.Content.Find.Execute("The software was created", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "The software was created", 2)

# ------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark from after "...diverse scales" to
#    right after "Software Description" in the title paragraph, and
#    merge the runs that used to straddle it into one run.
# ------------------------------------------------------------------

# Remove the existing bookmark (placed after "...diverse scales" in
# the commercial-applications paragraph).
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# Merge the two runs that used to be separated by the bookmark.
$d.Content.Find.Execute("scales, although", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "scales, although", 2)

# Re-create the bookmark at the end of the title paragraph
# ("Software Description"). A direct collapsed Range at the
# paragraph-mark position cannot host the bookmark reliably, so we
# temporarily append a one-character marker, bookmark that position
# (non-collapsed range), and then remove the marker text again —
# leaving the bookmark correctly anchored right after the run.
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertAfter("Z")

$marker = $d.Content
$marker.Find.Execute("Z")
$marker.Bookmarks.Add("_GoBack")
$marker.Text = ""
